# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 4
$ws.Range("Q2").Value = 1.77
$ws.Range("S2").Value = 2.88
$ws.Range("T2").Value = 1.73
$ws.Range("W2").Value = 2.14
$ws.Range("Y2").Value = 20
$ws.Range("AO2").Value = 46

# Row 3
$ws.Range("F3").Value = 3.55
$ws.Range("G3").Value = 5.1
$ws.Range("I3").Value = 2.1
$ws.Range("N3").Value = 3.9
$ws.Range("R3").Value = 1.4
$ws.Range("T3").Value = 1.62
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.9
$ws.Range("W3").Value = 1.24

# Row 4
$ws.Range("N4").Value = 2.66

# Row 6
$ws.Range("G6").Value = 2.88
$ws.Range("W6").Value = 1.53

# Row 7
$ws.Range("H7").Value = 1.54
$ws.Range("I7").Value = 1.65

# Row 8
$ws.Range("P8").Value = 1.9
$ws.Range("T8").Value = 2.06

# Row 10
$ws.Range("F10").Value = 2.92
$ws.Range("G10").Value = 2.98
$ws.Range("W10").Value = 1.5

# Row 11
$ws.Range("O11").Value = 1.32

# Row 12
$ws.Range("G12").Value = 2.26
$ws.Range("N12").Value = 2.28
$ws.Range("X12").Value = 6.4
$ws.Range("AB12").Value = 6
$ws.Range("AE12").Value = 110
$ws.Range("AG12").Value = 13
$ws.Range("AO12").Value = 180

# Row 13
$ws.Range("F13").Value = 1.88
$ws.Range("I13").Value = 5.3
$ws.Range("AD13").Value = 21
$ws.Range("AE13").Value = 75
